$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the trailing empty/border-formatting-only rows (originally rows 9-11)
$ws.Range("A9:F11").EntireRow.Delete()

# Remove the "locale" row (row 2: locale, EN, US, DE, FR, IT) - everything below shifts up
$ws.Rows.Item(2).Delete()
